# Rename "Delay In" / "Easing" options to "Ramp Time" / "Transition" across
# every per-effect category on the "Menu Mock" sheet, and refresh their
# tooltips to match.
#
#   OptionXxxDelayIn -> OptionXxxRampTime   (tooltip: duration of the ramp)
#   OptionXxxEasing   -> OptionXxxTransition (tooltip: curve-shape description)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

# Row -> new Option name (column B) and new Tooltip (column F).
$changes = @(
    @{ Row = 39; Option = "OptionBasicRampTime";      Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 40; Option = "OptionBasicTransition";     Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 48; Option = "OptionCriticalRampTime";    Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 49; Option = "OptionCriticalTransition";  Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 57; Option = "OptionDismemberRampTime";   Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 58; Option = "OptionDismemberTransition"; Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 66; Option = "OptionDecapRampTime";       Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 67; Option = "OptionDecapTransition";     Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 75; Option = "OptionLastEnemyRampTime";   Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 76; Option = "OptionLastEnemyTransition"; Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 83; Option = "OptionLastStandRampTime";   Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 84; Option = "OptionLastStandTransition";  Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" },
    @{ Row = 91; Option = "OptionParryRampTime";       Tooltip = "Duration of transition ramp (seconds)" },
    @{ Row = 92; Option = "OptionParryTransition";     Tooltip = "Curve shape for ramping into slow-mo (Off = instant)" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 2).Value = $change.Option   # column B = Option
    $ws.Cells.Item($change.Row, 6).Value = $change.Tooltip  # column F = Tooltip
}
